# Import the cleaned-up "Data" table into the active sheet (Taul1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Nimi"
$ws.Range("C1").Value = "Tsekkaus"

# --- Data rows --------------------------------------------------------------
$data = @(
    @(1, "Tonttu",     "Tämä jää"),
    @(2, "Toljander",  "Tämä jää"),
    @(3, "Joulupukki", "Tämä jää"),
    @(4, "Muori",      "Tämä jää"),
    @(5, "Smith",      "Tämä jää")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row = $row + 1
}

# --- Formatting -------------------------------------------------------------
# Columns A & B (ID, Nimi) + header for those columns: Calibri 15, automatic/black.
$ws.Range("A1:B6").Font.Size = 15

# Column C header ("Tsekkaus"): Calibri 15, red.
$ws.Range("C1").Font.Size = 15
$ws.Range("C1").Font.Color = 255

# Column C data ("Tämä jää" x5): Calibri 11, red.
$ws.Range("C2:C6").Font.Size = 11
$ws.Range("C2:C6").Font.Color = 255

# Row heights for all populated rows.
$ws.Range("A1:C6").RowHeight = 19.5

# --- Sheet view / selection --------------------------------------------------
$ws.Range("E3").Select() | Out-Null
